# "updated main GSC export data"
#
# The "Chart" sheet holds a rolling 90-day window of GSC export numbers
# (Date / Non-HTTPS URLs / Pages) in rows 2-91. Each refresh the window
# advances by one calendar day: the oldest day (2025-11-02) drops off the
# front, every remaining day's row shifts up by one, and the newest day
# (2026-01-31) takes over the last row's date label. The last row's
# "Pages" figure (column C) is left as-is -- no fresh value was appended
# for the new day in this export.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$firstRow = 2
$lastRow = 91
$newestDate = "2026-01-31"

# Snapshot the existing Date (col A) and Pages (col C) columns before
# writing anything, so the shift is computed from stable, pre-edit data.
# (.Value2 is used for reads -- a bare .Value getter doesn't resolve here.)
$dates = @{}
$pages = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $dates[$r] = $ws.Cells.Item($r, 1).Value2
    $pages[$r] = $ws.Cells.Item($r, 3).Value2
}

# Column A keeps storing plain text dates (as the export always has) --
# force Text format up front so Excel doesn't reinterpret the yyyy-MM-dd
# strings as date serials when they're written back.
$ws.Range("A$firstRow`:A$lastRow").NumberFormat = "@"

# Rows 2..90 take on the next row's old date and the next row's old
# "Pages" figure -- the sliding window shifting up by one row.
for ($r = $firstRow; $r -le ($lastRow - 1); $r++) {
    $ws.Cells.Item($r, 1).Value2 = $dates[$r + 1]
    $ws.Cells.Item($r, 3).Value2 = $pages[$r + 1]
}

# The last row becomes the newest day: advance its date label, keep its
# "Pages" value unchanged.
$ws.Cells.Item($lastRow, 1).Value2 = $newestDate
$ws.Cells.Item($lastRow, 3).Value2 = $pages[$lastRow]
